$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E3 (Lâmpada de Vapor Sódio - Potência Nominal): 15 -> 19
$ws.Range("E3").Value = 19

# Update D10 (Bomba Dosadora com Diafragma - Potência Nominal): 0.9 -> 0.45
$ws.Range("D10").Value = 0.45

# Update the selected range/active cell to F15
$ws.Range("F15").Select()
